# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$offSheet = $wb.Worksheets.Item("OFF")
$offSheet.Range("B2").Value = 495
$offSheet.Range("C2").Value = 350
$offSheet.Range("D2").Value = 122
$offSheet.Range("E2").Value = 49
$offSheet.Range("F2").Value = 7
$offSheet.Range("G2").Value = 5

# --- DEF sheet ---
$defSheet = $wb.Worksheets.Item("DEF")
$defSheet.Range("B2").Value = 490
$defSheet.Range("C2").Value = 321
$defSheet.Range("D2").Value = 69
$defSheet.Range("E2").Value = 29
$defSheet.Range("F2").Value = 6
$defSheet.Range("G2").Value = 5
